$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Alt1")
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Alt2"
Write-Host "Sheets:"
foreach ($s in $wb.Worksheets) { Write-Host " - " $s.Name }
Write-Host "B3 before copy:" $ws1.Range("B3").Value2
$ws1.Range("A1:N26").Copy($ws2.Range("A1"))
Write-Host "B3 after copy (Alt1):" $ws1.Range("B3").Value2
Write-Host "B3 after copy (Alt2):" $ws2.Range("B3").Value2
